# Add new progress as of date 04-Nov-2025
#   - Column H (PERIOD TO EXPIRE) decreases by 1 for every data row
#   - Column I (LAST UPDATE) moves from 03-Nov-2025 to 04-Nov-2025
#
# Column I holds the date as literal text (General-formatted inline string,
# not a real date serial). Assigning a date-shaped string straight to
# .Value makes Excel auto-detect it as a date (changing both the stored
# value and the cell's number format). To keep it as plain text we type it
# with a leading apostrophe (forces text), then immediately restore the
# cell's original formatting (which the apostrophe entry perturbs via the
# "quote prefix" flag) by copying the format from column C in the same
# row - a column that is never edited here, so it still carries the
# untouched style for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = 3..16

foreach ($r in $rows) {
    $hCell = $ws.Cells.Item($r, 8)   # column H
    $iCell = $ws.Cells.Item($r, 9)   # column I
    $fmtDonor = $ws.Cells.Item($r, 3)  # column C, same row - untouched style donor

    # PERIOD TO EXPIRE: one day less
    # (Range.Value's getter is unreliable in this host, so read/write via Value2)
    $hCell.Value2 = $hCell.Value2 - 1

    # LAST UPDATE: 03-Nov-2025 -> 04-Nov-2025, kept as literal text
    $iCell.Value2 = "'04-Nov-2025"

    $fmtDonor.Copy()
    $iCell.PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = $false
